$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Download Linux Mint Cinnamon Edition ISO" ->
#           "Download Linux Mint Cinnamon Edition ISOf" (typo "ISOf"),
#           split across three runs flanked by proofErr spell-check markers,
#           matching what Word's editor + spell-checker produce when a user
#           types an extra letter into a recognised word.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Download Linux Mint Cinnamon Edition ISO")
if ($found1) {
    $para1 = $r1.Paragraphs(1)
    $prange1 = $para1.Range

    $xmlFrag1 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="47752EFD" w14:textId="77777777" w:rsidR="00005BA2" w:rsidRDefault="0038163A"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Download Linux Mint Cinnamon Edition </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ISO</w:t></w:r><w:r><w:t>f</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    [void]$prange1.InsertXML($xmlFrag1)
}

# ---------------------------------------------------------------------------
# Change 2: merge the "- " and the marutter_pubkey URL runs into a single
#           run (same strike-through formatting, no text change).
# ---------------------------------------------------------------------------
$r2 = $d.Content
$old2 = "- https://cloud.r-project.org/bin/linux/ubuntu/marutter_pubkey.asc | "
[void]$r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# ---------------------------------------------------------------------------
# Change 3: merge the "apt-get -y " and "install make " runs into a single
#           run (same character style / font, no text change).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$old3 = "apt-get -y install make "
[void]$r3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

# ---------------------------------------------------------------------------
# Change 4: merge the three runs making up the very long oauth20_desktop.srf
#           redirect URL into a single run (no text change). The combined
#           string is far too long to hand to Find/Replace, so locate the
#           paragraph via a short, unique anchor, then delete + retype the
#           whole paragraph's text in one shot so it collapses into one run.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("https://login.live.com/oauth20_desktop.srf")
if ($found4) {
    $para4 = $r4.Paragraphs(1)
    $prange4 = $para4.Range
    $startPos4 = $prange4.Start

    $fullText4 = $prange4.Text
    # Trim the trailing paragraph-mark character(s) picked up by Paragraph.Range.
    $fullText4 = $fullText4.TrimEnd([char]13, [char]7)

    $delRange4 = $d.Range($startPos4, $startPos4 + $fullText4.Length)
    $delRange4.Delete()

    $insPoint4 = $d.Range($startPos4, $startPos4)
    [void]$insPoint4.InsertAfter($fullText4)
}
